$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Status")
Write-Output $ws.Name
Write-Output $ws.Range("A7").Value
